$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Sheet" (sheet1.xml) - flat data table, one row per object.
# Two duplicate rows are removed:
#   row 25 -> 105.6047547_0.8384163_10550_0020_24 (duplicate of row 7 / _6)
#   row 34 -> 105.4984167_0.0687500_10550_0020_33 (duplicate of row 8 / _7)
# After removal, every row below shifts up, and the trailing running
# index in column A ("..._N") is renumbered to stay sequential again.
# -----------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Sheet")

$wsData.Rows.Item(25).Delete()
$wsData.Rows.Item(33).Delete()

$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val) {
        $newVal = $val -replace '_\d+$', ('_' + ($r - 1))
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# -----------------------------------------------------------------
# Sheet "Sheet1" (sheet2.xml) - one labeled row every 13 rows, each
# block reserving room for a row of embedded thumbnail images drawn
# on top via the worksheet's drawing. Remove the same two duplicate
# entries by deleting their whole 13-row block, which keeps every
# other block's starting row (2, 15, 28, ...) fixed while shifting
# the block contents up - exactly like the canonical export does.
# -----------------------------------------------------------------
$wsImg = $wb.Worksheets.Item("Sheet1")

$wsImg.Rows("301:313").Delete()
$wsImg.Rows("405:417").Delete()

for ($r = 2; $r -le 600; $r += 13) {
    $cell = $wsImg.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val) {
        $n = (($r - 2) / 13) + 1
        $newVal = $val -replace '_\d+$', ('_' + $n)
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# Remove the 10 picture shapes (two groups of five filter images) that
# belonged to the two deleted catalog entries; they sat at the very
# end of the drawing (rows 612 and 625) and are not auto-removed by
# the row deletions above because they are one-cell-anchored.
$shapeCount = $wsImg.Shapes.Count
$stopAt = $shapeCount - 9
for ($i = $shapeCount; $i -ge $stopAt; $i--) {
    $wsImg.Shapes.Item($i).Delete()
}
